$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 870 (shifts old rows 870:911 down to 871:912)
$ws.Rows(870).Insert()

# Populate the newly inserted row with the new record (2026/02/25, 水, 8, 201).
# The date column is stored as plain text in this sheet (not a real date),
# so force text formatting before the write to stop Excel's autodetect from
# turning "2026/02/25" into a date serial, then restore the default "Normal"
# style so the cell matches its untouched siblings.
$ws.Range("A870").NumberFormat = "@"
$ws.Range("A870").Value = "2026/02/25"
$ws.Range("A870").Style = "Normal"

$ws.Range("B870").Value = "水"
$ws.Range("C870").Value = 8
$ws.Range("D870").Value = 201
